$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 76: B76 becomes a numeric 3 instead of text "3"; other cells stay the same
$ws.Range("B76").Value = 3

# Insert new row 77 with the new annotation data
$ws.Range("A77").Value = "Ruilin"
$ws.Range("B77").NumberFormat = "@"
$ws.Range("B77").Value = "3"
$ws.Range("B77").ClearFormats()
$ws.Range("C77").Value = "novel, hinder"
$ws.Range("D77").Value = "DIS"
$ws.Range("E77").Value = "OTH"
$ws.Range("F77").Value = "84d24e0a-0b18-4f4c-a441-4ea819712145"
$ws.Range("G77").Value = "rk9kKMZ0-_annotated.xlsx"
$ws.Range("H77").Value = "While the idea is novel and I do agree that I have not seen other works along these lines there are a few things that are missing and hinder this paper significantly."
